# Fix "Utgift" (expense) column E: the stored figures were accidentally
# multiplied by 12 (annualised) instead of being the monthly figures the
# rest of the model expects. Divide each by 12 on both the "private" and
# "Income" sheets, rows 7-72.
$wb = $excel.ActiveWorkbook

$sheetNames = @("private", "Income")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($r = 7; $r -le 72; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $old = $cell.Value2
        if ($old -ne $null) {
            $cell.Value2 = $old / 12
        }
    }

    # The shorter monthly figures need less column width than the old
    # annual ones did - mirror Excel's own best-fit recalculation for
    # column E (shrinks from ~16.4 to ~12.9 characters, matching the
    # already-"best fit" width used elsewhere in these sheets, e.g. col C).
    $ws.Columns.Item(5).ColumnWidth = 12
}
